$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6855646666666667
$ws.Range("H2").Value = 2.056694
$ws.Range("M2").Value = 0.3806726666666667
$ws.Range("N2").Value = 1.142018
$ws.Range("O2").Value = 0.1106000967880863
$ws.Range("P2").Value = 0.1106000967880863
$ws.Range("Q2").Value = 0.2609757298324444
$ws.Range("R2").Value = 2.348781568492
$ws.Range("S2").Value = 0.1106000967880863
$ws.Range("T2").Value = 0.1106000967880863

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6855646666666667
$ws.Range("H3").Value = 2.056694
$ws.Range("O3").Value = 0.04256506442063556
$ws.Range("P3").Value = 0.04256506442063556
$ws.Range("Q3").Value = 0.1004379659253333
$ws.Range("R3").Value = 0.9039416933280001
$ws.Range("S3").Value = 0.04256506442063556
$ws.Range("T3").Value = 0.04256506442063556

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6855646666666667
$ws.Range("H4").Value = 2.056694
$ws.Range("M4").Value = 1.310315666666667
$ws.Range("N4").Value = 3.930947
$ws.Range("O4").Value = 0.3806972557953004
$ws.Range("P4").Value = 0.3806972557953004
$ws.Range("Q4").Value = 0.8983061232464445
$ws.Range("R4").Value = 8.084755109218001
$ws.Range("S4").Value = 0.3806972557953004
$ws.Range("T4").Value = 0.3806972557953004

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6855646666666667
$ws.Range("H5").Value = 2.056694
$ws.Range("M5").Value = 1.604391333333333
$ws.Range("N5").Value = 4.813174
$ws.Range("O5").Value = 0.4661375829959777
$ws.Range("P5").Value = 0.4661375829959777
$ws.Range("Q5").Value = 1.099914009639556
$ws.Range("R5").Value = 9.899226086756002
$ws.Range("S5").Value = 0.4661375829959777
$ws.Range("T5").Value = 0.4661375829959777
